$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sequence Diagram")

$descriptions = @(
    "Creazione SD - Effettuare una ricerca",
    "Creazione SD - Aggiungere prodotto",
    "Creazione SD - Modificare quantita",
    "Creazione SD - Eliminare prodotto",
    "Creazione SD - Svuotare carrello",
    "Creazione SD - Effettuare ordine",
    "Creazione SD - Aggiornare dati"
)

# Source cell that already holds the text "0.2" (as a shared string) with the
# plain unformatted style used throughout the body of the table, so we can
# copy it down without Excel re-interpreting "0.2" as the number 0.2.
$versionSrc = $ws.Cells.Item(5, 3)

$startRow = 6
for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 2).Value = 44873
    $versionSrc.Copy($ws.Cells.Item($row, 3))
    $ws.Cells.Item($row, 4).Value = $descriptions[$i]
    $ws.Cells.Item($row, 5).Value = "SDS"
}

# The DESCRIZIONE column's per-row "best fit" style (s=4) that used to sit
# under the empty rows is no longer needed now that row 13 is the first
# blank row after the table body, so clear it back to the unformatted
# default, matching how Excel drops unused formatting here.
$ws.Range("D13:D23").Clear()

$ws.Range("D13").Select()
